$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "No. of R3 Excel Row's to Execute" value in D2 from 100 to 1
$ws.Range("D2").Value = "1"

# Move the selection / active cell to E24 as recorded in the saved view state
$ws.Range("E24").Select()
